$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31
$ws.Cells.Item($row, 1).Value = 45211
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 2).Value = "16:08"
$ws.Cells.Item($row, 3).Value = 76.7
$ws.Cells.Item($row, 4).Value = "natura"
